# Adding the "Assortment" KPI sheet to the CCAAU_SAND template workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the BayCountKPI "Value 2" cell: prepend the new scene type
#    "CCA 3 Door" to the existing comma separated list.
# ---------------------------------------------------------------------
$bayCountKPI = $wb.Worksheets.Item("BayCountKPI")
$bayCountKPI.Range("B2").Value = "CCA 3 Door, CCA - STANDARD LONG LANE CHECKOUT,CCA - EXPRESS CHECKOUT,CCA - SELF CHECK OUT,CCA - OTHER - FRONT OF STORE,CCA - INDIES - REST OF STORE,Competitor - STANDARD LONG LANE CHECKOUT,Competitor - EXPRESS CHECKOUT,Competitor - SELF CHECK OUT,Competitor - OTHER - FRONT OF STORE,Competitor - INDIES - REST OF STORE"
$bayCountKPI.Range("B2").Select()

# ---------------------------------------------------------------------
# 2. Create the new "Assortment" worksheet, positioned after BayCountKPI.
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Move($null, $wb.Worksheets.Item("BayCountKPI"))

# the sheet object becomes stale after Move(), so grab a fresh handle
$assortment = $wb.Worksheets.Item(4)
$assortment.Name = "Assortment"

# ---------------------------------------------------------------------
# 3. Header row formatting: copy the existing yellow / hairline-border
#    header style used on BayCountKPI!A1 so the style table is reused
#    instead of duplicated.
# ---------------------------------------------------------------------
$bayCountKPI.Range("A1").Copy()
$assortment.Range("A1:D1").PasteSpecial(-4122)
$assortment.Range("A1").VerticalAlignment = -4108
$assortment.Range("B1").HorizontalAlignment = -4108
$assortment.Range("B1").VerticalAlignment = -4108

# ---------------------------------------------------------------------
# 4. Header / data values.
# ---------------------------------------------------------------------
$assortment.Range("A1").Value = "scene_types_to_exclude"
$assortment.Range("B1").Value = "categories_to_exclude"
$assortment.Range("C1").Value = "brands_to_exclude"
$assortment.Range("D1").Value = "ean_codes_to_exclude"
$assortment.Range("A2").Value = "CCA-Standard checkout cooler, 2.2 Pepsi/Schweppes Impulse Cooler"

# ---------------------------------------------------------------------
# 5. Column widths to roughly match the authored layout.
# ---------------------------------------------------------------------
$assortment.Range("A1").ColumnWidth = 56.83
$assortment.Range("B1").ColumnWidth = 21.33
$assortment.Range("C1").ColumnWidth = 19.31
$assortment.Range("D1").ColumnWidth = 19.71

$assortment.Range("B2").Select()
